# Apply the commit's edits to the active workbook.
#
# Summary of the change (per the OOXML diff / commit message
# "Getting the right crossover now."):
#   1. The column-B header text (shared string used by B1) is corrected so
#      it actually matches the recurrence used by the B-column formulas:
#        "T(n) = 7T(n/2) + 2.25(n^2)"  ->  "T(n) = 7T(n/2)+18((n/4)^2)"
#   2. The sheet's saved cursor/selection moves from the old leftover
#      A15:XFD20 block selection to a single cell, H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the header label in B1 (shared string reused wherever B1 is referenced).
$ws.Range("B1").Value = "T(n) = 7T(n/2)+18((n/4)^2)"

# 2) Move the selection/active cell to H7.
$ws.Range("H7").Select()
